$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are preserved exactly as text
# (avoids Excel auto-converting "65.00" -> 65, "0.800" -> 0.8, etc.)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.934.11"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.503.69"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.19"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.49"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.41"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.944.66"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "58.860.07"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.78"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.494.25"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "322.58"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.94"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.00"
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.164"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.52"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "0.0₃0762"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.01"
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.46"
$ws.Range("E30").Value = "  -2.98%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.35"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.36"
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.52"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.58"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.800"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "281.40"
$ws.Range("E40").Value = "  +1.54%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.93"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.98"
$ws.Range("E44").Value = "  -5.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "128.73"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0499"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.24"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.749.83"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("E51").Value = "  -0.54%  "
